$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" column (G) values per regenerated save_data (K instead of Strike#)
$kValues = @{
    2 = 2
    3 = 0
    4 = 1
    5 = 0
    6 = 1
    7 = 1
    8 = 1
    9 = 1
    10 = 0
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 2
    25 = 3
    26 = 2
    27 = 1
    28 = 3
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 1
    34 = 0
    35 = 2
    36 = 1
    37 = 2
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 2
    44 = 1
    45 = 0
    46 = 2
    47 = 1
    48 = 0
    49 = 1
    50 = 3
    51 = 1
    52 = 1
    53 = 1
    54 = 3
    55 = 0
    56 = 2
    57 = 1
    58 = 0
    59 = 1
    60 = 0
    61 = 3
    62 = 0
    63 = 1
    64 = 1
    65 = 3
    66 = 2
    67 = 2
    68 = 2
    69 = 2
    71 = 1
    72 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Host "Updated $($kValues.Count) K column values"
